# Update the "Förändrad" (Changed) date column (C) for rows 2-12
# from 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244),
# matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = 45244
}
